# Apply the "only one descriptor/value pair per line is copied to
# description files" edit to INSTRUCTIONS.docx.
#
# Most of the underlying change in the source commit is a re-wrap of
# paragraph runs (splitting long runs of plain text into several runs
# with identical text/formatting - a no-op for the rendered document).
# The only textual changes that affect what the reader actually sees
# are:
#   1) a stray comma turned into a period after
#      "...the study ID (Short name)..."
#   2) the sentence about the prime symbol, which now reads
#      "...it is safer to spell it, like 5-prime." instead of
#      "...use ' instead."
#
# We locate each passage with Find (exact text, no wildcards) and then
# assign the replacement directly to the found Range's .Text property
# (instead of using Find's built-in Replace) so that Word's
# smart-quote autocorrect does not mangle the straight apostrophes
# that appear in the surrounding text.

$d = $word.ActiveDocument

function Replace-ExactText($oldText, $newText) {
    $range = $d.Content
    $found = $range.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find text: $oldText"
    }
    $range.Text = $newText
}

# 1) "and enter the study ID (Short name), This file ..." -> "...(Short name). This file ..."
Replace-ExactText `
    "and enter the study ID (Short name), This file will make a new directory tree, rooted in the" `
    "and enter the study ID (Short name). This file will make a new directory tree, rooted in the"

# 2) Prime-symbol sentence rewritten.
Replace-ExactText `
    ") are allowed. Be careful if the description contains prime symbol (' ,as in 5'): use ' instead." `
    ") are allowed. Be careful if the description contains prime symbol (' ,as in 5'), it is safer to spell it, like 5-prime."
